$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "24/09/2021" entries become real dates, note gets a "V 0/0 -" prefix ---
# A6 already carries the date-formatted style; just replace its free-text date with a real date value.
$ws.Cells.Item(6, 1).Value2 = 44463

# B6 becomes a date too; give it A6's date number format by copying formats across.
$ws.Cells.Item(6, 2).Value2 = 44463
$ws.Range("A6").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats

# C6 keeps its existing note, just prefixed with "V 0/0 - "
$ws.Cells.Item(6, 3).Value2 = "V 0/0 - " + $ws.Cells.Item(6, 3).Value2

# D6 (the long repair description) is untouched.

# --- Row 7: "30/09/2021" entries become real dates, note gets a "V 0/0 -" prefix ---
$ws.Cells.Item(7, 1).Value2 = 44469
$ws.Cells.Item(7, 2).Value2 = 44469

# C7 keeps its existing note, just prefixed with "V 0/0 - "
$ws.Cells.Item(7, 3).Value2 = "V 0/0 - " + $ws.Cells.Item(7, 3).Value2

# D7 (the long repair description) is untouched.

# --- Restore the view: selection moved on to C7 (scrolls frozen pane up to A6) ---
$ws.Activate()
$ws.Range("C7").Select()
